# Fruta / hortaliza, semanal
# A new weekly price record is prepended as row 3, pushing the existing
# rows 3-16 down to rows 4-17 (data content unchanged, just relocated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3; this shifts old rows 3..16 down to 4..17
# and inherits formatting (e.g. the date style) from the row being pushed down.
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with this week's new record.
$ws.Cells.Item(3, 1).Value = 11
$ws.Cells.Item(3, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(3, 3).Value = "Bíobío"
$ws.Cells.Item(3, 4).Value = 44532
$ws.Cells.Item(3, 5).Value = 8
$ws.Cells.Item(3, 6).Value = 100112022
$ws.Cells.Item(3, 7).Value = "Arveja Verde"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 250
$ws.Cells.Item(3, 11).Value = 14000
$ws.Cells.Item(3, 12).Value = 15000
$ws.Cells.Item(3, 13).Value = 14400
$ws.Cells.Item(3, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(3, 15).Value = "Región del Maule"
$ws.Cells.Item(3, 16).Value = 576
$ws.Cells.Item(3, 17).Value = 25
$ws.Cells.Item(3, 18).Value = "Hortaliza"
